# Apply the edit described by the diff:
# Insert two new data rows (new weekly observations) right before the
# existing row 401, shifting the old rows 401-425 down to 403-427.
# The two new rows hold "Camote" price observations for origin "Peru"
# dated 44516 (1a nueva(o) / 2a nueva(o)).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 401 (existing content shifts down).
$ws.Rows.Item(401).EntireRow.Insert(-4121)
$ws.Rows.Item(401).EntireRow.Insert(-4121)

# --- New row 401 ---
$ws.Cells.Item(401, 1).Value = 8
$ws.Cells.Item(401, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(401, 3).Value = "Coquimbo"
$ws.Cells.Item(401, 4).Value = 44516
$ws.Cells.Item(401, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(401, 5).Value = 4
$ws.Cells.Item(401, 6).Value = 100112045
$ws.Cells.Item(401, 7).Value = "Zapallo"
$ws.Cells.Item(401, 8).Value = "Camote"
$ws.Cells.Item(401, 9).Value = "1a nueva(o)"
$ws.Cells.Item(401, 10).Value = 800
$ws.Cells.Item(401, 11).Value = 700
$ws.Cells.Item(401, 12).Value = 750
$ws.Cells.Item(401, 13).Value = 725
$ws.Cells.Item(401, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(401, 15).Value = "Perú"
$ws.Cells.Item(401, 16).Value = 725
$ws.Cells.Item(401, 17).Value = 1
$ws.Cells.Item(401, 18).Value = "Hortaliza"

# --- New row 402 ---
$ws.Cells.Item(402, 1).Value = 8
$ws.Cells.Item(402, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(402, 3).Value = "Coquimbo"
$ws.Cells.Item(402, 4).Value = 44516
$ws.Cells.Item(402, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(402, 5).Value = 4
$ws.Cells.Item(402, 6).Value = 100112045
$ws.Cells.Item(402, 7).Value = "Zapallo"
$ws.Cells.Item(402, 8).Value = "Camote"
$ws.Cells.Item(402, 9).Value = "2a nueva(o)"
$ws.Cells.Item(402, 10).Value = 520
$ws.Cells.Item(402, 11).Value = 600
$ws.Cells.Item(402, 12).Value = 650
$ws.Cells.Item(402, 13).Value = 625
$ws.Cells.Item(402, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(402, 15).Value = "Perú"
$ws.Cells.Item(402, 16).Value = 625
$ws.Cells.Item(402, 17).Value = 1
$ws.Cells.Item(402, 18).Value = "Hortaliza"

Write-Host ("Final used range: " + $ws.UsedRange.Address())
